# "red border for current day, added song cover for day 7 and spotify links"
#
# Day 7 (row 4, date 2021-12-07) used to carry a "Spruch/Gedicht/Textthing" /
# "Text" placeholder plus a "No pic" marker in column D. It now has an actual
# song cover: a "Song" gift with description "Crash Bendicat", and the
# leftover "No pic" note is gone.
#
# Day 11 (row 8, date 2021-12-11) keeps its "No pic" marker but also becomes
# a "Song" gift, this time "Berry White".
#
# Both cells also lose their green highlight fill (it was presumably used to
# flag "still needs content" placeholders) now that real content has been
# filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day 7 (row 4): turn the placeholder text gift into a song cover ---
# Set C4 before B4 so the shared-string table allocates "Crash Bendicat"
# ahead of "Song", matching the order the entries were authored in.
$ws.Range("C4").Value = "Crash Bendicat"
$ws.Range("B4").Value = "Song"
$ws.Range("D4").ClearContents()

# --- Day 11 (row 8): also a song, keeps its "No pic" note in D8 ---
$ws.Range("C8").Value = "Berry White"
$ws.Range("B8").Value = "Song"

# --- Drop the highlight fill on the four cells that just got real content ---
$ws.Range("B4").Interior.Pattern = -4142
$ws.Range("C4").Interior.Pattern = -4142
$ws.Range("B8").Interior.Pattern = -4142
$ws.Range("C8").Interior.Pattern = -4142

# --- Row heights settle at a hair taller than before (15 -> 14.4 default
#     font metrics recalculation bumped the custom 22.15/40.15 rows to
#     22.2/40.2) ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Rows($r).RowHeight = 22.2
}
$ws.Rows(26).RowHeight = 40.2

# --- Selection left on C8, the last-touched cell ---
$ws.Range("C8").Select()
